$d = $word.ActiveDocument

function Replace-ParagraphRunText($doc, $oldText, $newText) {
    # Find every top-level paragraph whose text (ignoring the trailing
    # paragraph mark) equals $oldText exactly, and swap that text for
    # $newText while preserving any sibling runs (e.g. leading empty
    # <w:r/> runs) and any *direct* (explicit) bold/italic character
    # formatting applied to the run that held the text - formatting
    # that merely comes from the paragraph style (e.g. Heading1 being
    # bold) must NOT be re-applied as direct formatting.
    foreach ($p in $doc.Paragraphs) {
        $pText = $p.Range.Text
        $pText = $pText.TrimEnd([char]13, [char]7)
        if ($pText -eq $oldText) {
            $start = $p.Range.Start
            $end = $p.Range.End - 1
            $r = $doc.Range($start, $end)

            $runBold = $r.Font.Bold
            $runItalic = $r.Font.Italic

            $styleFont = $p.Style.Font
            $styleBold = $styleFont.Bold
            $styleItalic = $styleFont.Italic

            $directBold = ($runBold -and -not $styleBold)
            $directItalic = ($runItalic -and -not $styleItalic)

            $r.InsertXML($newText)

            $r2 = $doc.Range($start, $start + $newText.Length)
            if ($directBold) { $r2.Font.Bold = $runBold }
            if ($directItalic) { $r2.Font.Italic = $runItalic }
        }
    }
}

Replace-ParagraphRunText $d "Play Big Bad Wolf for Free - A Fairytale Themed Slot Game" "Play Big Bad Wolf Free and Enjoy Immersive Fairytale Gameplay"

Replace-ParagraphRunText $d "Fairytale theme based on The Three Little Pigs story" "Exceptional graphics and music"

Replace-ParagraphRunText $d "Immersive graphics and music" "Immersive fairytale theme"

Replace-ParagraphRunText $d "Standard 5x3 configuration and 25 pay lines" "Variety of gameplay features"

Replace-ParagraphRunText $d "Wild symbols, Scatter symbols, and Bonus feature available" "Free Spins and Bonus feature for additional winnings"

Replace-ParagraphRunText $d "Lack of progressive jackpot" "Limited number of pay lines"

Replace-ParagraphRunText $d "Free Spins feature does not offer a large number of spins" "May not appeal to players who prefer non-fairytale themes"

Replace-ParagraphRunText $d "Read our review of Big Bad Wolf, an online slot game developed by Quickspin and inspired by The Three Little Pigs story. Play this fairy tale-themed game for free now!" "Read our review of Big Bad Wolf slot game and play for free. Experience the immersive fairytale theme."
